$wb = $excel.ActiveWorkbook

# New technologies to append to both the withdrawals and consumption
# "Electricity Source" subscript sheets (issues #280 and #99).
# Each entry is: new label text, and the existing row (on the same sheet)
# whose value the new row's formula should reuse.
$newRows = @(
    @{ Label = "hard coal w CCS";                     RefRow = 2  },
    @{ Label = "natural gas combined cycle w CCS";     RefRow = 4  },
    @{ Label = "biomass w CCS";                        RefRow = 10 },
    @{ Label = "lignite w CCS";                         RefRow = 14 },
    @{ Label = "small modular reactor";                 RefRow = 5  },
    @{ Label = "hydrogen";                               RefRow = 4  }
)

$sheetNames = @("WUbPPT-withdrawals", "WUbPPT-consumption")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $startRow = 19
    for ($i = 0; $i -lt $newRows.Count; $i++) {
        $row = $startRow + $i
        $entry = $newRows[$i]

        $ws.Cells.Item($row, 1).Value = $entry.Label
        $ws.Cells.Item($row, 2).Formula = "=B" + $entry.RefRow
    }

    $ws.Range("A25").Select()
}
